$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4").Value = "03 a 05"
$ws.Range("A5").Value = "11 a 20"
$ws.Range("A6").Value = "06 a 10"
$ws.Range("A7").Value = "50 o más"
$ws.Range("A8").Value = "21 a 50"
